# This script reassigns the species/find-specific data (Id, Taxonsorteringsordning,
# TaxonId, Artnamn, Vetenskapligt namn, Auktor, Alder-Stadium, Kon, Aktivitet, Metod,
# Ost, Nord, Publik kommentar, Substrat-beskrivning) across rows 2-17 of the active
# sheet, while leaving the location/visit metadata columns (C, D, I, P, S, T, U, V, W,
# Y, Z, AA, AB, AD, AE, AG, AI, AT, AW, AX, AY) untouched for each row.
#
# All values written below are literal, taken from the *original* (pre-edit) state of
# the corresponding source row, so the writes can be applied in any order without a
# "clobber before read" hazard.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 4
$ws.Range("A2").Value = 111438439
$ws.Range("B2").Value = 78107
$ws.Range("E2").Value = 6453
$ws.Range("F2").Value = 'Vedskivlav'
$ws.Range("G2").Value = 'Hertelidea botryosa'
$ws.Range("H2").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("Q2").Value = 468788.4775288465
$ws.Range("R2").Value = 6882785.67140964
$ws.Range("AC2").ClearContents()
$ws.Range("AO2").Value = 'silverved tall'

# Row 3 <- original row 7
$ws.Range("A3").Value = 111438428
$ws.Range("B3").Value = 77597
$ws.Range("E3").Value = 864
$ws.Range("F3").Value = 'Knottrig blåslav'
$ws.Range("G3").Value = 'Hypogymnia bitteri'
$ws.Range("H3").Value = '(Lynge) Ahti'
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("Q3").Value = 468740.5586073888
$ws.Range("R3").Value = 6882780.957796668
$ws.Range("AC3").ClearContents()
$ws.Range("AO3").Value = 'björk'

# Row 4 <- original row 11
$ws.Range("A4").Value = 111438432
$ws.Range("B4").Value = 77267
$ws.Range("E4").Value = 6446
$ws.Range("F4").Value = 'Kolflarnlav'
$ws.Range("G4").Value = 'Carbonicola anthracophila'
$ws.Range("H4").Value = '(Nyl.) Bendiksby & Timdal'
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("Q4").Value = 468756.5460031229
$ws.Range("R4").Value = 6882784.091042386
$ws.Range("AC4").ClearContents()
$ws.Range("AO4").Value = 'brandstubbe'

# Row 5 <- original row 17
$ws.Range("A5").Value = 111438426
$ws.Range("B5").Value = 76918
$ws.Range("E5").Value = 6437
$ws.Range("F5").Value = 'Blanksvart spiklav'
$ws.Range("G5").Value = 'Calicium denigratum'
$ws.Range("H5").Value = '(Vain.) Tibell'
$ws.Range("K5").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("Q5").Value = 468629.2461709682
$ws.Range("R5").Value = 6882722.464435354
$ws.Range("AC5").ClearContents()
$ws.Range("AO5").Value = 'silverved tall'

# Row 6 <- original row 15
$ws.Range("A6").Value = 111438453
$ws.Range("B6").Value = 78107
$ws.Range("E6").Value = 6453
$ws.Range("F6").Value = 'Vedskivlav'
$ws.Range("G6").Value = 'Hertelidea botryosa'
$ws.Range("H6").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("K6").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()
$ws.Range("Q6").Value = 468789.3971357156
$ws.Range("R6").Value = 6882885.489071017
$ws.Range("AC6").ClearContents()
$ws.Range("AO6").Value = 'silverved tall'

# Row 7 <- original row 16
$ws.Range("A7").Value = 111438433
$ws.Range("B7").Value = 78081
$ws.Range("E7").Value = 229821
$ws.Range("F7").Value = 'Vedflamlav'
$ws.Range("G7").Value = 'Ramboldia elabens'
$ws.Range("H7").Value = '(Fr.) Kantvilas & Elix'
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("Q7").Value = 468756.5460031229
$ws.Range("R7").Value = 6882784.091042386
$ws.Range("AC7").ClearContents()
$ws.Range("AO7").Value = 'silverved tall'

# Row 8 <- original row 10
$ws.Range("A8").Value = 111438430
$ws.Range("B8").Value = 77268
$ws.Range("E8").Value = 228912
$ws.Range("F8").Value = 'Mörk kolflarnlav'
$ws.Range("G8").Value = 'Carbonicola myrmecina'
$ws.Range("H8").Value = '(Ach.) Bendiksby & Timdal'
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("Q8").Value = 468756.5460031229
$ws.Range("R8").Value = 6882784.091042386
$ws.Range("AC8").ClearContents()
$ws.Range("AO8").Value = 'brandstubbe'

# Row 9 <- original row 8
$ws.Range("A9").Value = 111438442
$ws.Range("B9").Value = 73696
$ws.Range("E9").Value = 6440
$ws.Range("F9").Value = 'Vitgrynig nållav'
$ws.Range("G9").Value = 'Chaenotheca subroscida'
$ws.Range("H9").Value = '(Eitner) Zahlbr.'
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("Q9").Value = 468800.3617588138
$ws.Range("R9").Value = 6882801.965499061
$ws.Range("AC9").ClearContents()
$ws.Range("AO9").Value = 'gran'

# Row 10 <- original row 2
$ws.Range("A10").Value = 111438447
$ws.Range("B10").Value = 76495
$ws.Range("E10").Value = 6487
$ws.Range("F10").Value = 'Blågrå svartspik'
$ws.Range("G10").Value = 'Chaenothecopsis fennica'
$ws.Range("H10").Value = '(Laurila) Tibell'
$ws.Range("K10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("Q10").Value = 468866.1318338988
$ws.Range("R10").Value = 6882808.390505624
$ws.Range("AC10").ClearContents()
$ws.Range("AO10").Value = 'silverved tall'

# Row 11 <- original row 14
$ws.Range("A11").Value = 111438455
$ws.Range("B11").Value = 77515
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = 'Garnlav'
$ws.Range("G11").Value = 'Alectoria sarmentosa'
$ws.Range("H11").Value = '(Ach.) Ach.'
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("Q11").Value = 468784.2260541836
$ws.Range("R11").Value = 6882884.599394682
$ws.Range("AC11").ClearContents()
$ws.Range("AO11").Value = 'tall'

# Row 12 <- original row 9
$ws.Range("A12").Value = 111438440
$ws.Range("B12").Value = 73689
$ws.Range("E12").Value = 308
$ws.Range("F12").Value = 'Brunpudrad nållav'
$ws.Range("G12").Value = 'Chaenotheca gracillima'
$ws.Range("H12").Value = '(Vain.) Tibell'
$ws.Range("K12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("Q12").Value = 468800.2970216064
$ws.Range("R12").Value = 6882794.936009536
$ws.Range("AC12").ClearContents()
$ws.Range("AO12").Value = 'högstubbe björk'

# Row 14 <- original row 5
$ws.Range("A14").Value = 111438446
$ws.Range("B14").Value = 77550
$ws.Range("E14").Value = 185
$ws.Range("F14").Value = 'Violettgrå tagellav'
$ws.Range("G14").Value = 'Bryoria nadvornikiana'
$ws.Range("H14").Value = '(Gyeln.) Brodo & D.Hawksw.'
$ws.Range("K14").ClearContents()
$ws.Range("L14").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("Q14").Value = 468853.3954624244
$ws.Range("R14").Value = 6882801.477506777
$ws.Range("AC14").ClearContents()
$ws.Range("AO14").Value = 'gran'

# Row 15 <- original row 6
$ws.Range("A15").Value = 111438425
$ws.Range("B15").Value = 56398
$ws.Range("E15").Value = 100109
$ws.Range("F15").Value = 'Tretåig hackspett'
$ws.Range("G15").Value = 'Picoides tridactylus'
$ws.Range("H15").Value = '(Linnaeus, 1758)'
$ws.Range("K15").Value = ''
$ws.Range("L15").Value = ''
$ws.Range("M15").Value = 'färska spår'
$ws.Range("N15").Value = ''
$ws.Range("Q15").Value = 468571.5178632676
$ws.Range("R15").Value = 6882722.999468728
$ws.Range("AC15").Value = 'Ringhack i tall'
$ws.Range("AO15").Value = 'tall'

# Row 16 <- original row 3
$ws.Range("A16").Value = 111438444
$ws.Range("B16").Value = 77515
$ws.Range("E16").Value = 6425
$ws.Range("F16").Value = 'Garnlav'
$ws.Range("G16").Value = 'Alectoria sarmentosa'
$ws.Range("H16").Value = '(Ach.) Ach.'
$ws.Range("K16").ClearContents()
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("Q16").Value = 468841.2361184616
$ws.Range("R16").Value = 6882806.276033297
$ws.Range("AC16").ClearContents()
$ws.Range("AO16").Value = 'tall'

# Row 17 <- original row 12
$ws.Range("A17").Value = 111438457
$ws.Range("B17").Value = 78107
$ws.Range("E17").Value = 6453
$ws.Range("F17").Value = 'Vedskivlav'
$ws.Range("G17").Value = 'Hertelidea botryosa'
$ws.Range("H17").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("K17").ClearContents()
$ws.Range("L17").ClearContents()
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("Q17").Value = 468747.5763832342
$ws.Range("R17").Value = 6882880.250689426
$ws.Range("AC17").ClearContents()
$ws.Range("AO17").Value = 'silverved tall'

